$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Locate the paragraph that immediately precedes the bookmark-only
#    paragraph (the one holding the "_GoBack" bookmark).  We find it by its
#    exact text rather than a hard-coded index, so the script stays correct
#    even if paragraph numbering shifts a little.
# ---------------------------------------------------------------------------
$anchorText = "На сторінці вмінь змінено шаблон відображення, перший варіант дизайну з градієнтним розділювачем між зображенням та інформативним блоком."

$count = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    $tt = $t.TrimEnd([char]13, [char]7)
    if ($tt -eq $anchorText) {
        $anchorIndex = $i
    }
}

$bookmarkParaIndex = $anchorIndex + 1
$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIndex)
$bookmarkRange = $bookmarkPara.Range

# ---------------------------------------------------------------------------
# 2. Build the OOXML for the four new changelog paragraphs.  The final
#    paragraph in this fragment has no trailing paragraph mark of its own in
#    the inserted body, so InsertXML merges its runs into the following
#    (pre-existing) paragraph -- the one that carries the "_GoBack" bookmark
#    -- leaving the bookmark paragraph itself intact with the new runs
#    placed in front of the bookmark.
# ---------------------------------------------------------------------------
$newParagraphsXml =
    '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Зролено</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> сторінку про себе.</w:t></w:r></w:p>' +
    '<w:p><w:r><w:t>Виправлено стилі елементів що збились.</w:t></w:r></w:p>' +
    '<w:p><w:r><w:t>Перевірена адаптивність сторінки вмінь</w:t></w:r></w:p>' +
    '<w:p><w:r><w:t xml:space="preserve">Змінено </w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>header</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>адаптивн</w:t></w:r>' +
    '<w:r><w:t>ість</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>, відкорегованостилі.</w:t></w:r></w:p>'

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
    $newParagraphsXml +
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint = $d.Range($bookmarkRange.Start, $bookmarkRange.Start)
$null = $insertionPoint.InsertXML($packageXml)

# ---------------------------------------------------------------------------
# 3. Insert one extra, genuinely empty paragraph right after the bookmark
#    paragraph (so it ends up between that paragraph and the document's
#    pre-existing trailing empty paragraph).  InsertXML collapses a lone
#    empty "<w:p/>" fragment, so two are inserted and the extra one is then
#    deleted again -- this leaves a clean "<w:p/>" with no stray run.
# ---------------------------------------------------------------------------
$bookmarkParaIndex = $anchorIndex + 4
$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIndex)
$afterBookmark = $d.Range($bookmarkPara.Range.End, $bookmarkPara.Range.End)

$twoEmptyParasXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$null = $afterBookmark.InsertXML($twoEmptyParasXml)

$extraParaIndex = $bookmarkParaIndex + 2
$extraPara = $d.Paragraphs.Item($extraParaIndex)
$null = $extraPara.Range.Delete()

Write-Host "Done"
